# Edit: update Junction Flooding data table (custom accuracy + 1000 new data rows commit)
# - Overwrite rows 2-5 with new sensor readings
# - Remove the now-obsolete row 6 (dimension shrinks to A1:AH5)
# - Widen several data columns by 1 character (col width 7->8, col T 8->9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45164.50694444445, 10.726, 7.333, 3.404, 23.56, 17.15, 8.176, 24.228, 13.347, 5.245, 7.323, 9.309, 10.191, 2.44, 8.647, 11.655, 7.955, 2.648, 1.093, 124.223, 23.834, 7.982, 14.964, 8.049, 2.19, 13.597, 7.05, 6.629, 7.562, 9.942, 2.682, 21.628, 4.075, 9.978),
    @(45164.51388888889, 22.362, 16.53, 1.958, 49.032, 39.31, 17.476, 65.143, 27.31, 12.014, 17.42, 19.636, 20.996, 5.455, 17.67, 24.933, 15.174, 1.388, 0.951, 261.642, 49.368, 16.31, 32.835, 17.275, 2.726, 32.958, 14.407, 12.913, 15.128, 20.666, 1.136, 59.418, 9.025, 20.392),
    @(45164.52083333334, 19.534, 14.526, 1.409, 42.814, 34.562, 15.297, 61.092, 23.821, 10.58, 15.357, 17.163, 18.324, 4.787, 15.414, 21.824, 13.155, 0.963, 0.743, 227.314, 43.14, 14.228, 28.809, 15.127, 2.274, 29.905, 12.568, 11.219, 13.164, 18.051, 0.722, 55.614, 7.919, 17.789),
    @(45164.52777777778, 23.89, 17.85, 1.35, 52.24, 42.56, 18.75, 73.32, 29.06, 12.99, 19.02, 20.95, 22.3, 5.91, 18.8, 26.71, 15.87, 0.8, 0.82, 278.82, 52.57, 17.35, 35.31, 18.54, 2.65, 35.91, 15.33, 13.59, 15.97, 22.03, 0.53, 66.59, 9.75, 21.69),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}

# Row 6 no longer exists in the updated dataset - delete it (shifts dimension to A1:AH5)
$ws.Rows.Item(6).Delete()

# Column width tweaks (ColumnWidth setter stores width + ~0.8333 padding in the
# saved XML, so subtract that offset to land exactly on the target character width)
$ws.Columns.Item(2).ColumnWidth = 7.165
$ws.Columns.Item(3).ColumnWidth = 7.165
$ws.Columns.Item(5).ColumnWidth = 7.165
$ws.Columns.Item(6).ColumnWidth = 7.165
$ws.Columns.Item(7).ColumnWidth = 7.165
$ws.Columns.Item(9).ColumnWidth = 7.165
$ws.Columns.Item(10).ColumnWidth = 7.165
$ws.Columns.Item(11).ColumnWidth = 7.165
$ws.Columns.Item(12).ColumnWidth = 7.165
$ws.Columns.Item(13).ColumnWidth = 7.165
$ws.Columns.Item(15).ColumnWidth = 7.165
$ws.Columns.Item(16).ColumnWidth = 7.165
$ws.Columns.Item(17).ColumnWidth = 7.165
$ws.Columns.Item(20).ColumnWidth = 8.165
$ws.Columns.Item(21).ColumnWidth = 7.165
$ws.Columns.Item(22).ColumnWidth = 7.165
$ws.Columns.Item(23).ColumnWidth = 7.165
$ws.Columns.Item(24).ColumnWidth = 7.165
$ws.Columns.Item(26).ColumnWidth = 7.165
$ws.Columns.Item(27).ColumnWidth = 7.165
$ws.Columns.Item(28).ColumnWidth = 7.165
$ws.Columns.Item(29).ColumnWidth = 7.165
$ws.Columns.Item(30).ColumnWidth = 7.165
$ws.Columns.Item(34).ColumnWidth = 7.165
